$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(701,  3, 90, 45, 97, 15),
    @(601,  9, 60, 67, 60, 42),
    @(201,  9, 30, 15, 45, 30),
    @(901, 16, 15, 45, 60, 60),
    @(501,  9, 52, 30, 75, 45),
    @(801,  3, 67, 65, 52, 45),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6, 45, 30, 60, 45),
    @(401,  9, 48, 67, 75, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(101,  9, 30, 15, 60, 15),
    @(1001,18, 30, 75, 60, 72),
    @(1101, 0, 15, 30, 30,  0),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(1,    0,  2,  2,  2,  2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
